# Scheduled-runner refresh: push updated Universalis price snapshots
# (currentAveragePrice / NQ / HQ) and the recomputed Leve profit columns
# for the affected leves across each crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1405.44
$ws.Range("I98").Value = 1739.1875
$ws.Range("K98").Value = 1739.1875
$ws.Range("M98").Value = -241.1875

$ws.Range("H122").Value = 1405.44
$ws.Range("I122").Value = 1739.1875
$ws.Range("K122").Value = 5217.5625
$ws.Range("M122").Value = -2767.5625

$ws.Range("H132").Value = 5850292.5
$ws.Range("I132").Value = 7248698
$ws.Range("K132").Value = 21746094
$ws.Range("M132").Value = -21743564

$ws.Range("H137").Value = 1484.0358
$ws.Range("I137").Value = 1243
$ws.Range("J137").Value = 2930.25
$ws.Range("K137").Value = 3729
$ws.Range("L137").Value = 8790.75
$ws.Range("M137").Value = -1179
$ws.Range("N137").Value = -13890.75

$ws.Range("H138").Value = 1725.1959
$ws.Range("I138").Value = 683.05884
$ws.Range("J138").Value = 1946.65
$ws.Range("K138").Value = 2049.17652
$ws.Range("L138").Value = 5839.950000000001
$ws.Range("M138").Value = 3090.82348
$ws.Range("N138").Value = -16119.95

$ws.Range("H139").Value = 43445.715
$ws.Range("J139").Value = 43445.715
$ws.Range("L139").Value = 43445.715
$ws.Range("N139").Value = -53725.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1165.1428
$ws.Range("I61").Value = 982.5263
$ws.Range("K61").Value = 982.5263
$ws.Range("M61").Value = -770.5263

$ws.Range("H110").Value = 1358.3
$ws.Range("I110").Value = 1165.9474
$ws.Range("K110").Value = 1165.9474
$ws.Range("M110").Value = 879.0526

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws.Range("H132").Value = 2389.8838
$ws.Range("I132").Value = 2128.4595
$ws.Range("J132").Value = 4002
$ws.Range("K132").Value = 6385.3785
$ws.Range("L132").Value = 12006
$ws.Range("M132").Value = -3855.3785
$ws.Range("N132").Value = -17066

$ws.Range("H136").Value = 1165.1428
$ws.Range("I136").Value = 982.5263
$ws.Range("K136").Value = 2947.5789
$ws.Range("M136").Value = -397.5789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 99999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 299997
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -305067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1054.875

$ws.Range("H28").Value = 29643
$ws.Range("J28").Value = 29643
$ws.Range("L28").Value = 29643
$ws.Range("N28").Value = -30133

$ws.Range("H31").Value = 1215.6428
$ws.Range("I31").Value = 866.913
$ws.Range("K31").Value = 866.913
$ws.Range("M31").Value = -571.913

$ws.Range("H34").Value = 1215.6428
$ws.Range("I34").Value = 866.913
$ws.Range("K34").Value = 866.913
$ws.Range("M34").Value = -664.913

$ws.Range("H95").Value = 26900
$ws.Range("J95").Value = 26900
$ws.Range("L95").Value = 26900
$ws.Range("N95").Value = -32392

$ws.Range("H132").Value = 12927.556
$ws.Range("I132").Value = 12927.556
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 38782.66800000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -36252.66800000001
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 12346807
$ws.Range("I134").Value = 12821588
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 38464764
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -38462229
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1476.8334
$ws.Range("I5").Value = 2051
$ws.Range("J5").Value = 673
$ws.Range("K5").Value = 6153
$ws.Range("L5").Value = 2019
$ws.Range("M5").Value = -6041
$ws.Range("N5").Value = -2243

$ws.Range("H17").Value = 954.5454999999999
$ws.Range("I17").Value = 833.3333
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 2499.9999
$ws.Range("L17").Value = 4500
$ws.Range("M17").Value = -2330.9999
$ws.Range("N17").Value = -4838

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H99").Value = 2046.5834
$ws.Range("I99").Value = 625
$ws.Range("J99").Value = 2520.4443
$ws.Range("K99").Value = 1875
$ws.Range("L99").Value = 7561.3329
$ws.Range("M99").Value = 371
$ws.Range("N99").Value = -12053.3329

$ws.Range("H135").Value = 1476.8334
$ws.Range("I135").Value = 2051
$ws.Range("J135").Value = 673
$ws.Range("K135").Value = 18459
$ws.Range("L135").Value = 6057
$ws.Range("M135").Value = -15924
$ws.Range("N135").Value = -11127

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H111").Value = 50020000
$ws.Range("I111").Value = 100000000
$ws.Range("K111").Value = 100000000
$ws.Range("M111").Value = -99996933

$ws.Range("H132").Value = 2768.45
$ws.Range("I132").Value = 2298.389
$ws.Range("K132").Value = 6895.167
$ws.Range("M132").Value = -4365.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2413.2
$ws.Range("J2").Value = 2301.2856
$ws.Range("L2").Value = 2301.2856
$ws.Range("N2").Value = -2525.2856

$ws.Range("H40").Value = 2596
$ws.Range("I40").Value = 2596
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2596
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2460
$ws.Range("N40").ClearContents()

$ws.Range("H122").Value = 28343742
$ws.Range("I122").Value = 28343742
$ws.Range("K122").Value = 85031226
$ws.Range("M122").Value = -85028776

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 21503750
$ws.Range("J5").Value = 21503750
$ws.Range("L5").Value = 21503750
$ws.Range("N5").Value = -21503974

$ws.Range("H123").Value = 56252.332
$ws.Range("J123").Value = 56252.332
$ws.Range("L123").Value = 56252.332
$ws.Range("N123").Value = -66052.33199999999

$ws.Range("H126").Value = 50506480
$ws.Range("I126").Value = 123457600
$ws.Range("K126").Value = 370372800
$ws.Range("M126").Value = -370370330

$ws.Range("H132").Value = 2014.7727
$ws.Range("I132").Value = 1821.1714
$ws.Range("J132").Value = 2767.6667
$ws.Range("K132").Value = 5463.5142
$ws.Range("L132").Value = 8303.000100000001
$ws.Range("M132").Value = -2933.5142
$ws.Range("N132").Value = -13363.0001

$ws.Range("H136").Value = 595.1053000000001
$ws.Range("I136").Value = 311.47827
$ws.Range("J136").Value = 1030
$ws.Range("K136").Value = 934.43481
$ws.Range("L136").Value = 3090
$ws.Range("M136").Value = 1615.56519
$ws.Range("N136").Value = -8190

